$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "31.100.40"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "1.951.76"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'245.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.4874"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'44.65"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.2961"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'0.06810"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'19.02"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").Value = "'106.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.73%  "
$ws.Range("D13").Value = "1.931.95"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'0.07720"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "'0.7109"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.46%  "
$ws.Range("D17").Value = "'286.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.46%  "
$ws.Range("D18").Value = "30.986.36"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "'0.000007738"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'13.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.189.44"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'5.508"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'6.589"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'9.896"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "'168.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'19.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'2.186"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").Value = "'0.1050"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "'1.438"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").Value = "'4.715"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +14.86%  "
$ws.Range("D33").Value = "'4.455"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("D34").Value = "'0.04996"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "'0.7601"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'1.158"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'0.02037"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'2.702"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'2.142"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.30%  "
$ws.Range("D41").Value = "'6.404"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.71%  "
$ws.Range("D42").Value = "'0.4469"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "'109.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'0.8791"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "'72.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").Value = "'0.9986"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'7.448"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'979.05"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +15.64%  "
$ws.Range("D49").Value = "'0.1277"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("D50").Value = "'9.368"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'0.2584"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.53%  "
